$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column C (shifts old C..G to D..H)
$ws.Columns("C").Insert()

# Fill in the new row 4 (new user story / task)
$ws.Range("A4").Value = "I want to be notified when creating an account was successful"
$ws.Range("B4").Value = "Modify server side and GUI for account creation success"

# New column header (merge C1:C2 like the other header cells)
$ws.Range("C1:C2").Merge()
$ws.Range("C1").Value = "Task Assigned"

# Fill remaining new user stories (column A)
$ws.Range("A5").Value = "I want to share my pictures to the public"
$ws.Range("A6").Value = "I want to share my pictures privately to some users"
$ws.Range("A7").Value = "I want to delete images from server side"

# Fill remaining new task descriptions (column B)
$ws.Range("B5").Value = "Implement publicly shared functionality"
$ws.Range("B6").Value = "Implement privately shared functionality"
$ws.Range("B7").Value = "Implement delete functionality for server side"

# Update the burndown chart series reference to follow the shifted columns
$chart = $ws.ChartObjects(1).Chart
$chart.SeriesCollection(1).Formula = "=SERIES(,,Sheet1!`$D`$27:`$H`$27,1)"

# Shift the chart's anchor position right by one column (it was anchored relative
# to the old column layout, and does not automatically follow the inserted column)
$co = $ws.ChartObjects(1)
$co.Left = $co.Left + $ws.Columns("C").Width
